# Bill of Material update:
#  - remove the two "Need to order" placeholder rows (JST connector, old battery)
#  - keep the rows that follow (5-way switch / prototyping parts / SHIP) which
#    shift up into their place
#  - record the actual order: 4 Li-po batteries (2000mAh) and 5 right-angle
#    JST-PH connectors, plus a shipping line item

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "JST-PH connector" and "Lithium battery (1200mAh)" rows that
# were still marked "Need to order" - row 10 twice since the second delete
# pulls the next row up into row 10's place.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

# Row 13: shipping line item
$ws.Range("A13").Value = "Shipping"
$ws.Range("B13").Value = "mixec"
$ws.Range("C13").Value = 9.15
$ws.Range("D13").Value = "USPS"
$ws.Range("E13").Value = "Nathaniel"
$ws.Range("F13").Value = "Batteries and connectors from Adafruit"

# Row 14: the new battery, 4 ordered
$ws.Range("A14").Value = "Lithium Ion Polymer Battery - 3.7v 2000mAh"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 50
$ws.Range("D14").Value = "Adafruit"
$ws.Range("E14").Value = "Nathaniel"

# Row 15: the JST-PH right-angle connector, 5 ordered
$ws.Range("A15").Value = "JST-PH 2-Pin SMT Right Angle Connector"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 3.75
$ws.Range("D15").Value = "Adafruit"
$ws.Range("E15").Value = "Nathaniel"
$ws.Range("F15").Value = "Battery connector which mounts on the PCB."
